$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'263.44"

$ws.Range("D3").Value = "'22.79"

$ws.Range("D4").Value = "'6.204"

$ws.Range("D5").Value = "'0.06097"

$ws.Range("D6").Value = "'3.516"

$ws.Range("D7").Value = "'6.712"

$ws.Range("D8").Value = "'1.363"

$ws.Range("D9").Value = "'0.7988"

$ws.Range("D11").Value = "'0.08148"

$ws.Range("D12").Value = "'0.03321"

$ws.Range("D14").Value = "'0.09251"

$ws.Range("D15").Value = "'3.895"

$ws.Range("D16").Value = "'0.001690"

$ws.Range("D17").Value = "'0.04826"

$ws.Range("D18").Value = "'0.0006208"

$ws.Range("D19").Value = "'0.006191"

$ws.Range("D20").Value = "'0.005990"

$ws.Range("D22").Value = "'0.0001502"

$ws.Range("D23").Value = "'3.693"

$ws.Range("D24").Value = "'2.279"

$ws.Range("D25").Value = "'0.3386"

$ws.Range("D26").Value = "'0.1268"

$ws.Range("D27").Value = "'0.0003247"

$ws.Range("D40").Value = "'0.04652"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1120"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003134"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.007279"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").Value = "'0.01095"

$ws.Range("B45").Value = "ACDXExchange"
$ws.Range("C45").Value = "https://coinranking.com/coin/-y35lbZ7U+acdxexchange-acxt"
$ws.Range("D45").Value = "'0.002974"
$ws.Range("E45").Value = "44ACDXExchangeACXTBestin24h"

$ws.Range("B46").Value = "CoinLion"
$ws.Range("C46").Value = "https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion"
$ws.Range("D46").Value = "'0.00006053"
$ws.Range("E46").Value = "45CoinLionLION"

$ws.Range("B47").Value = "Kangarootoken"
$ws.Range("C47").Value = "https://coinranking.com/coin/zkVNkSGwZ3+kangarootoken-gar"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "46KangarootokenGAR"

$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.7009"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOIN"

$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.1541"
$ws.Range("E49").Value = "48BOLOBOLO"

$ws.Range("B50").Value = "CryptobidCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "49CryptobidCoinCBC"

$ws.Range("B51").Value = "SpecialPowerGold"
$ws.Range("C51").Value = "https://coinranking.com/coin/jPTWzmsWb+specialpowergold-spg"
$ws.Range("D51").Value = "'0.01011"
$ws.Range("E51").Value = "50SpecialPowerGoldSPG"
